# Apply "hybrid bold + color highlighting" to quantitative metrics in the
# achievement / responsibility bullet paragraphs, per the commit diff
# ("Implement quantitative metrics highlighting across all resume formats").
#
# For each target paragraph we know (from the diff) the exact ordered
# sequence of text chunks and whether each chunk should be rendered
# bold + colored (RGB 2C3E50) or left as plain text. We locate the
# paragraph by its original full text (robust against index drift), then
# overwrite it chunk-by-chunk: the first chunk replaces the whole
# paragraph via Range.Text, and each subsequent chunk is appended with
# Range.InsertAfter and immediately formatted. This reproduces the
# run-split structure shown in the diff (alternating plain /
# <w:b/><w:color w:val="2C3E50"/> runs).

function Get-WdColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$MetricColor = Get-WdColor "2C3E50"

function Get-ParagraphText($p) {
    # Paragraph.Range.Text includes the trailing paragraph-mark character;
    # strip it so callers can compare against plain text.
    $t = $p.Range.Text
    if ($t.Length -gt 0) {
        return $t.Substring(0, $t.Length - 1)
    }
    return $t
}

function Find-ParagraphIndex($d, [string]$originalText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ((Get-ParagraphText $p) -eq $originalText) {
            return $i
        }
    }
    throw "Could not locate paragraph with text: $originalText"
}

function Set-ParagraphRuns($d, [string]$originalText, [object[]]$chunks) {
    # $chunks is an array of hashtables: @{ Text = "..."; Bold = $true/$false }
    # Sanity check: the concatenation of chunk text must equal the paragraph's
    # current text (minus the trailing paragraph mark), so we never silently
    # touch the wrong paragraph or drop/alter characters.
    $expected = ""
    foreach ($c in $chunks) { $expected += $c.Text }
    if ($expected -ne $originalText) {
        throw "Chunk reconstruction mismatch. Expected [$originalText] got [$expected]"
    }

    $paraIndex = Find-ParagraphIndex $d $originalText

    # First chunk: replace the whole paragraph (Range includes the trailing
    # paragraph mark, but setting .Text on it keeps exactly one paragraph).
    $first = $chunks[0]
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.Text = $first.Text
    if ($first.Bold) {
        $p = $d.Paragraphs.Item($paraIndex)
        $chunkStart = $p.Range.Start
        $chunkEnd = $chunkStart + $first.Text.Length
        $r = $d.Range($chunkStart, $chunkEnd)
        $r.Font.Bold = 1
        $r.Font.Color = $MetricColor
    }

    # Remaining chunks: insert right before the paragraph mark, then format.
    for ($i = 1; $i -lt $chunks.Count; $i++) {
        $chunk = $chunks[$i]
        $p = $d.Paragraphs.Item($paraIndex)
        $insertAt = $p.Range.End - 1
        $ip = $d.Range($insertAt, $insertAt)
        $ip.InsertAfter($chunk.Text)
        if ($chunk.Bold) {
            $chunkEnd = $insertAt + $chunk.Text.Length
            $r = $d.Range($insertAt, $chunkEnd)
            $r.Font.Bold = 1
            $r.Font.Color = $MetricColor
        }
    }
}

$d = $word.ActiveDocument

Set-ParagraphRuns $d `
    '• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%' `
    @(
        @{ Text = '• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from '; Bold = $false },
        @{ Text = '23%'; Bold = $true },
        @{ Text = ' to '; Bold = $false },
        @{ Text = '64%'; Bold = $true }
    )

Set-ParagraphRuns $d `
    '• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes' `
    @(
        @{ Text = '• Utilized advanced sampling methods to decrease survey margin of error from '; Bold = $false },
        @{ Text = '±4.2%'; Bold = $true },
        @{ Text = ' to '; Bold = $false },
        @{ Text = '±2.1%'; Bold = $true },
        @{ Text = ', increasing voter turnout prediction accuracy from '; Bold = $false },
        @{ Text = '71%'; Bold = $true },
        @{ Text = ' to '; Bold = $false },
        @{ Text = '87%'; Bold = $true },
        @{ Text = ', and ensuring survey results more closely reflected true population attitudes'; Bold = $false }
    )

Set-ParagraphRuns $d `
    '• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis' `
    @(
        @{ Text = '• Trigonometric algorithm for boundary estimation reduced mapping costs by '; Bold = $false },
        @{ Text = '73.5%'; Bold = $true },
        @{ Text = ', saving campaigns and organizations '; Bold = $false },
        @{ Text = '$4.7M'; Bold = $true },
        @{ Text = ' and enabling smaller nonprofits to conduct analysis'; Bold = $false }
    )

Set-ParagraphRuns $d `
    '• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion' `
    @(
        @{ Text = '• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over '; Bold = $false },
        @{ Text = '$2'; Bold = $true },
        @{ Text = ' trillion'; Bold = $false }
    )

Set-ParagraphRuns $d `
    '• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%' `
    @(
        @{ Text = '• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by '; Bold = $false },
        @{ Text = '57%'; Bold = $true }
    )

Set-ParagraphRuns $d `
    '• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations' `
    @(
        @{ Text = '• Platform impact: Built redistricting system serving '; Bold = $false },
        @{ Text = '12,847'; Bold = $true },
        @{ Text = ' analysts across 89 organizations'; Bold = $false }
    )

Set-ParagraphRuns $d `
    '• Revenue generation: Delivered $4.9M additional revenue through optimization' `
    @(
        @{ Text = '• Revenue generation: Delivered '; Bold = $false },
        @{ Text = '$4.9M'; Bold = $true },
        @{ Text = ' additional revenue through optimization'; Bold = $false }
    )

Set-ParagraphRuns $d `
    '• 23% conversion rate improvement' `
    @(
        @{ Text = '• '; Bold = $false },
        @{ Text = '23%'; Bold = $true },
        @{ Text = ' conversion rate improvement'; Bold = $false }
    )

Write-Output "Applied quantitative metrics highlighting to 8 bullet paragraphs."
